# Applies the "Updated symbol list" data refresh (Thu Jan 5 09:10:23 UTC 2023).
# Source sheet stores every data cell as literal text (prices/volumes/hour look
# numeric but are inline strings), so new values are entered with a leading
# apostrophe to force text entry instead of Excel's automatic number/percent/date
# coercion; the quote-prefix formatting that introduces is cleared at the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'256.77"
$ws.Range("E2").Value = "'0.53%"
$ws.Range("G2").Value = "'9"

# Row 3: OKB
$ws.Range("D3").Value = "'27.11"
$ws.Range("E3").Value = "'-3.38%"
$ws.Range("G3").Value = "'9"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'4.754"
$ws.Range("E4").Value = "'-10.05%"
$ws.Range("G4").Value = "'9"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.05944"
$ws.Range("E5").Value = "'1.52%"
$ws.Range("G5").Value = "'9"

# Row 6: KuCoinToken
$ws.Range("D6").Value = "'6.656"
$ws.Range("E6").Value = "'-0.69%"
$ws.Range("G6").Value = "'9"

# Row 7: MXToken
$ws.Range("D7").Value = "'0.8705"
$ws.Range("E7").Value = "'0.27%"
$ws.Range("G7").Value = "'9"

# Row 8: FTXToken
$ws.Range("D8").Value = "'0.9471"
$ws.Range("E8").Value = "'-0.19%"
$ws.Range("G8").Value = "'9"

# Row 9: WazirX
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-0.57%"
$ws.Range("G9").Value = "'9"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.03942"
$ws.Range("E10").Value = "'14.07%"
$ws.Range("G10").Value = "'9"

# Row 11: MandalaExchangeToken
$ws.Range("D11").Value = "'0.07166"
$ws.Range("E11").Value = "'0.75%"
$ws.Range("G11").Value = "'9"

# Row 12: BitrueCoin
$ws.Range("D12").Value = "'0.03198"
$ws.Range("E12").Value = "'0.66%"
$ws.Range("G12").Value = "'9"

# Row 13: BitMartToken
$ws.Range("D13").Value = "'0.09262"
$ws.Range("E13").Value = "'0.37%"
$ws.Range("G13").Value = "'9"

# Row 14: BitForexToken
$ws.Range("D14").Value = "'0.001552"
$ws.Range("E14").Value = "'-0.36%"
$ws.Range("G14").Value = "'9"

# Row 15: TigerCash
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006054"
$ws.Range("E15").Value = "'4.29%"
$ws.Range("G15").Value = "'9"

# Row 16: LEO
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.482"
$ws.Range("E16").Value = "'-0.41%"
$ws.Range("G16").Value = "'9"

# Row 17: GateToken
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.199"
$ws.Range("E17").Value = "'-1.04%"
$ws.Range("G17").Value = "'9"

# Row 18: BTSEToken
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.219"
$ws.Range("E18").Value = "'0.73%"
$ws.Range("G18").Value = "'9"

# Row 19: One
$ws.Range("B19").Value = "One"
$ws.Range("C19").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D19").Value = "'0.01047"
$ws.Range("E19").Value = "'1,631.08%"
$ws.Range("G19").Value = "'9"

# Row 20: BitpandaEcosystemToken
$ws.Range("E20").Value = "'-1.40%"
$ws.Range("G20").Value = "'9"

# Row 21: ProBitToken
$ws.Range("E21").Value = "'-0.77%"
$ws.Range("G21").Value = "'9"

# Row 22: MCDex
$ws.Range("D22").Value = "'3.804"
$ws.Range("E22").Value = "'7.58%"
$ws.Range("G22").Value = "'9"

# Row 23: CoinExToken
$ws.Range("D23").Value = "'0.04203"
$ws.Range("E23").Value = "'1.36%"
$ws.Range("G23").Value = "'9"

# Row 24: ZBToken
$ws.Range("E24").Value = "'2.61%"
$ws.Range("G24").Value = "'9"

# Row 25: BitKan
$ws.Range("D25").Value = "'0.001221"
$ws.Range("E25").Value = "'-1.03%"
$ws.Range("G25").Value = "'9"

# Row 26: HotbitToken
$ws.Range("D26").Value = "'0.004493"
$ws.Range("E26").Value = "'-8.57%"
$ws.Range("G26").Value = "'9"

# Row 27: NitroEx
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("E27").Value = "'0.05%"
$ws.Range("G27").Value = "'9"

# Row 28: UpBots
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("E28").Value = "'142.24%"
$ws.Range("G28").Value = "'9"

# Row 29: Spectre.aiUtilityToken
$ws.Range("G29").Value = "'9"

# Row 30: LegolasExchange
$ws.Range("G30").Value = "'9"

# Row 31: BitZToken
$ws.Range("G31").Value = "'9"

# Row 32: Birake
$ws.Range("G32").Value = "'9"

# Row 33: NashExchange
$ws.Range("G33").Value = "'9"

# Row 34: AAXToken
$ws.Range("G34").Value = "'9"

# Row 35: CenX
$ws.Range("G35").Value = "'9"

# Row 36: BNIXToken
$ws.Range("G36").Value = "'9"

# Row 37: Polkally
$ws.Range("G37").Value = "'9"

# Row 38: Charli3
$ws.Range("G38").Value = "'9"

# Row 39: BlubitexToken
$ws.Range("G39").Value = "'9"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.03832"
$ws.Range("E40").Value = "'0.32%"
$ws.Range("G40").Value = "'9"

# Row 41: KickToken
$ws.Range("D41").Value = "'0.006113"
$ws.Range("E41").Value = "'16.64%"
$ws.Range("G41").Value = "'9"

# Row 42: BKEXToken
$ws.Range("E42").Value = "'-0.14%"
$ws.Range("G42").Value = "'9"

# Row 43: CEJI
$ws.Range("D43").Value = "'0.002253"
$ws.Range("E43").Value = "'-3.67%"
$ws.Range("G43").Value = "'9"

# Row 44: LocalTraders
$ws.Range("D44").Value = "'0.01056"
$ws.Range("E44").Value = "'4.27%"
$ws.Range("G44").Value = "'9"

# Row 45: CoinLion
$ws.Range("D45").Value = "'0.00005502"
$ws.Range("E45").Value = "'5.33%"
$ws.Range("G45").Value = "'9"

# Row 46: Kangarootoken
$ws.Range("E46").Value = "'0.04%"
$ws.Range("G46").Value = "'9"

# Row 47: CoinbaseStockToken
$ws.Range("D47").Value = "'0.08856"
$ws.Range("E47").Value = "'-4.77%"
$ws.Range("G47").Value = "'9"

# Row 48: BOLO
$ws.Range("D48").Value = "'0.002389"
$ws.Range("E48").Value = "'10.95%"
$ws.Range("G48").Value = "'9"

# Row 49: CryptobidCoin
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("G49").Value = "'9"

# Row 50: SpecialPowerGold
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("G50").Value = "'9"

# Row 51: DigiFinexToken
$ws.Range("G51").Value = "'9"

# Clear the quote-prefix style stamped onto the cells above so formatting stays
# identical to the untouched header (row 1) and index column (A).
$ws.Range("B2:G51").Style = "Normal"
